{"js": "// Apply the commit's text changes:\n//   1. Title paragraph: \"Earthquakes\" -> \"Earthquakes- mod Carol\"\n//   2. First Author paragraph: \"Steve Purves\" -> \"Steve Purves adpted Musso\"\n//\n// The target OOXML keeps the existing runs untouched except the final\n// \"Earthquakes\" run (whose text grows a trailing \"-\") and then appends\n// four new sibling runs (\" \", \"mod\", \" \", \"Carol\") / (\" \", \"adpted\", \" \",\n// \"Musso\") to each paragraph, so we mirror that with a single in-place\n// insertText(\"Replace\") on the run's own range, followed by paragraph-end\n// insertText calls (each call mints its own new run).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Title paragraph: \"La Palma Earthquakes\" ------------------------\nconst titlePara = paragraphs.items[0];\n\n// Find the \"Earthquakes\" run inside the title paragraph only (scope the\n// search to that paragraph's own range so other \"Earthquakes\" elsewhere\n// in the document are left untouched).\nconst titleRange = titlePara.getRange();\nconst titleHits = titleRange.search(\"Earthquakes\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\ntitleHits.items[0].insertText(\"Earthquakes-\", \"Replace\");\nawait context.sync();\n\ntitlePara.insertText(\" \", \"End\");\ntitlePara.insertText(\"mod\", \"End\");\ntitlePara.insertText(\" \", \"End\");\ntitlePara.insertText(\"Carol\", \"End\");\nawait context.sync();\n\n// --- 2. First Author paragraph: \"Steve Purves\" --------------------------\nconst authorPara = paragraphs.items[1];\n\nauthorPara.insertText(\" \", \"End\");\nauthorPara.insertText(\"adpted\", \"End\");\nauthorPara.insertText(\" \", \"End\");\nauthorPara.insertText(\"Musso\", \"End\");\nawait context.sync();\n", "ps1": "# Apply the commit's text changes:\n#   1. Title paragraph:  \"La Palma Earthquakes\" -> \"La Palma Earthquakes- mod Carol\"\n#   2. First Author paragraph: \"Steve Purves\" -> \"Steve Purves adpted Musso\"\n#\n# The target OOXML keeps every existing run untouched except the\n# \"Earthquakes\" run (text grows a trailing \"-\") and then appends four new\n# sibling runs (\" \", \"mod\", \" \", \"Carol\") / (\" \", \"adpted\", \" \", \"Musso\")\n# to each paragraph. Word's `Range.Text = ...` setter in this host merges\n# the edited span together with whatever follows it in the paragraph into\n# a single run, so instead we delete the old word's Range and\n# `InsertAfter` the replacement at the now-collapsed caret -- that mints a\n# clean, separate run exactly like the recorded edit, without disturbing\n# the neighboring \" \" / \"Palma\" runs. New trailing tokens are appended the\n# same way, one call per run, right before the paragraph mark.\n\n$d = $word.ActiveDocument\n\nfunction Append-Run($para, [string]$token) {\n    # Insert immediately before the paragraph mark so the new text stays\n    # inside this paragraph (End sits just after the mark).\n    $e = $para.Range.End\n    $caret = $d.Range($e - 1, $e - 1)\n    $caret.InsertAfter($token)\n}\n\n# --- 1. Title paragraph --------------------------------------------------\n$titlePara = $d.Paragraphs(1)\n\n# Locate the \"Earthquakes\" run, scoped to this paragraph only, via Find on\n# a duplicate range (so the stored $titlePara.Range is left alone).\n$seek = $titlePara.Range.Duplicate\n$seek.Find.ClearFormatting()\n$found = $seek.Find.Execute(\"Earthquakes\", $true)\nif (-not $found) {\n    throw \"Could not find 'Earthquakes' in the title paragraph\"\n}\n\n# Replace just that run's text: delete the old word, then insert the new\n# text at the collapsed caret left behind (keeps the surrounding \" \" and\n# \"Palma\" runs untouched, unlike Range.Text= which merges forward).\n$wordStart = $seek.Start\n$old = $d.Range($wordStart, $seek.End)\n$old.Delete()\n$caret = $d.Range($wordStart, $wordStart)\n$caret.InsertAfter(\"Earthquakes-\")\n\nAppend-Run $titlePara \" \"\nAppend-Run $titlePara \"mod\"\nAppend-Run $titlePara \" \"\nAppend-Run $titlePara \"Carol\"\n\n# --- 2. First Author paragraph (\"Steve Purves\") --------------------------\n$authorPara = $d.Paragraphs(2)\n\nAppend-Run $authorPara \" \"\nAppend-Run $authorPara \"adpted\"\nAppend-Run $authorPara \" \"\nAppend-Run $authorPara \"Musso\"\n"}
